$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AP3").Value = 1.87
$ws.Range("AQ3").Value = 2.03
$ws.Range("G4").Value = 3
$ws.Range("I4").Value = 2.8
$ws.Range("J4").Value = 4
$ws.Range("AA4").Value = 13
$ws.Range("AB4").Value = 34
$ws.Range("AH4").Value = 101
$ws.Range("AK4").Value = 11
$ws.Range("AN4").Value = 29
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 3.75
$ws.Range("J5").Value = 2.25
$ws.Range("K5").Value = 2.05
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("S5").Value = 4.33
$ws.Range("T5").Value = 1.2
$ws.Range("U5").Value = 1.5
$ws.Range("V5").Value = 2.5
$ws.Range("W5").Value = 2.25
$ws.Range("X5").Value = 1.57
$ws.Range("AB5").Value = 11
$ws.Range("AC5").Value = 15
$ws.Range("AE5").Value = 7.5
$ws.Range("AF5").Value = 7.5
$ws.Range("AH5").Value = 81
$ws.Range("AJ5").Value = 12
$ws.Range("AL5").Value = 19
$ws.Range("AO5").Value = 51
$ws.Range("AP5").Value = 1.78
$ws.Range("AQ5").Value = 2.1
$ws.Range("AR5").Value = 3.65
$ws.Range("AS5").Value = 1.29
$ws.Range("K6").Value = 1.8
$ws.Range("AA6").Value = 11
$ws.Range("AF6").Value = 6
$ws.Range("AK6").Value = 17
$ws.Range("AR6").Value = 5.2
$ws.Range("AS6").Value = 1.16
$ws.Range("N7").Value = 4.75
$ws.Range("Y7").Value = 5
$ws.Range("AD7").Value = 51
$ws.Range("AE7").Value = 4.75
$ws.Range("Q8").Value = 2.88
$ws.Range("R8").Value = 1.4
$ws.Range("G11").Value = 1.5
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 6.5
$ws.Range("J11").Value = 2.05
$ws.Range("L11").Value = 6
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("W11").Value = 1.83
$ws.Range("X11").Value = 1.83
$ws.Range("AG11").Value = 17
$ws.Range("AI11").Value = 301
$ws.Range("AM11").Value = 67
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.9
$ws.Range("Q27").Value = 1.98
$ws.Range("R27").Value = 1.88
$ws.Range("O30").Value = 1.33
$ws.Range("P30").Value = 3.25
$ws.Range("Q30").Value = 2.05
$ws.Range("R30").Value = 1.75
